$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 752-753, pushing the previous rows 752.. down to 754..
$ws.Rows.Item(752).Resize(2).Insert()

# New row 752: Navel Late / Primera
$ws.Cells.Item(752, 1).Value = 9
$ws.Cells.Item(752, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(752, 3).Value = "Metropolitana"
$ws.Cells.Item(752, 4).Value = 44769
$ws.Cells.Item(752, 5).Value = 13
$ws.Cells.Item(752, 6).Value = "Fruta"
$ws.Cells.Item(752, 7).Value = 100102
$ws.Cells.Item(752, 8).Value = "Cítricos"
$ws.Cells.Item(752, 9).Value = 100102005
$ws.Cells.Item(752, 10).Value = "Naranja"
$ws.Cells.Item(752, 11).Value = "Navel Late"
$ws.Cells.Item(752, 12).Value = "Primera"
$ws.Cells.Item(752, 13).Value = 80
$ws.Cells.Item(752, 14).Value = 6000
$ws.Cells.Item(752, 15).Value = 6000
$ws.Cells.Item(752, 16).Value = 6000
$ws.Cells.Item(752, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(752, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(752, 19).Value = 333
$ws.Cells.Item(752, 20).Value = 18

# New row 753: Navel Late / Segunda
$ws.Cells.Item(753, 1).Value = 9
$ws.Cells.Item(753, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(753, 3).Value = "Metropolitana"
$ws.Cells.Item(753, 4).Value = 44769
$ws.Cells.Item(753, 5).Value = 13
$ws.Cells.Item(753, 6).Value = "Fruta"
$ws.Cells.Item(753, 7).Value = 100102
$ws.Cells.Item(753, 8).Value = "Cítricos"
$ws.Cells.Item(753, 9).Value = 100102005
$ws.Cells.Item(753, 10).Value = "Naranja"
$ws.Cells.Item(753, 11).Value = "Navel Late"
$ws.Cells.Item(753, 12).Value = "Segunda"
$ws.Cells.Item(753, 13).Value = 100
$ws.Cells.Item(753, 14).Value = 5000
$ws.Cells.Item(753, 15).Value = 5000
$ws.Cells.Item(753, 16).Value = 5000
$ws.Cells.Item(753, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(753, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(753, 19).Value = 278
$ws.Cells.Item(753, 20).Value = 18

# Ensure the date number format is applied like the other D-column cells
$ws.Cells.Item(752, 4).NumberFormat = $ws.Cells.Item(754, 4).NumberFormat
$ws.Cells.Item(753, 4).NumberFormat = $ws.Cells.Item(754, 4).NumberFormat
